# "Generate Report for Handback"
#
# This applies the handback-report update to localization-status.xlsx:
#   1. Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#      on the Overview sheet and on each locale sheet (zh-cn, de-de).
#   2. The "Latest Handback DateTime" (column H) on the locale sheets is
#      stamped with the real handback time (previously the zero-date
#      placeholder 0001-01-01 00:00:00).
#   3. The new "Latest Target File" (F) / "Latest Handback File" (G)
#      columns are populated with hyperlinked file names for each row,
#      mirroring the existing "Source File Name" (A) / "Latest Handoff
#      File" (D) hyperlinks for that row.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

# Latest Handback DateTime now has a real timestamp.
$zh.Range("H2").Value = "2016-03-19 16:45:56"
$zh.Range("H3").Value = "2016-03-19 16:45:56"

# Latest Target File (F) / Latest Handback File (G) - same files as the
# Source File Name (A) / Latest Handoff File (D) hyperlinks for the row.
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/15a3b9792f88db9c5e8f4ee5150a039059ba1bb6/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md", "", "", "18633056-ff34-44e2-8461-f8cb6b01ebaa.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c87ab28f4397f23bdecdace7cdc91d00fd45d36b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf", "", "", "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/15a3b9792f88db9c5e8f4ee5150a039059ba1bb6/e2e/6100965e-3277-4e74-8ceb-b89abe4613f0.md", "", "", "6100965e-3277-4e74-8ceb-b89abe4613f0.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c87ab28f4397f23bdecdace7cdc91d00fd45d36b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.zh-cn.xlf", "", "", "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# Latest Handback DateTime now has a real timestamp (de-de handback
# happened a little later than zh-cn).
$de.Range("H2").Value = "2016-03-19 16:46:11"
$de.Range("H3").Value = "2016-03-19 16:46:11"

# Latest Target File (F) / Latest Handback File (G).
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/15a3b9792f88db9c5e8f4ee5150a039059ba1bb6/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md", "", "", "18633056-ff34-44e2-8461-f8cb6b01ebaa.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/877f5e547c7e74a229a60ade298dc1c5a504a455/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf", "", "", "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/15a3b9792f88db9c5e8f4ee5150a039059ba1bb6/e2e/6100965e-3277-4e74-8ceb-b89abe4613f0.md", "", "", "6100965e-3277-4e74-8ceb-b89abe4613f0.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/877f5e547c7e74a229a60ade298dc1c5a504a455/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.de-de.xlf", "", "", "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.de-de.xlf") | Out-Null

Write-Output "Handback report generated."
